# Refresh the "cryptos" price/volume table (GitHub Actions crypto-price bot run).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). All four are text cells, not
# numbers, so numeric-looking prices are written with a leading "'" (forces
# text entry, preserving trailing zeros like "39.70") and then the cell style
# is reset to "Normal" so no stray number-format style gets attached.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '49.599.03'
$ws.Cells.Item(2, 5).Value = '  -0.73%  '
$ws.Cells.Item(3, 4).Value = '2.635.89'
$ws.Cells.Item(3, 5).Value = '  -0.57%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).Value = "'112.28"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.82%  '
$ws.Cells.Item(6, 4).Value = "'324.35"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.89%  '
$ws.Cells.Item(7, 4).Value = "'0.525"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.11%  '
$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).Value = "'0.545"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.62%  '
$ws.Cells.Item(10, 4).Value = "'39.70"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.09%  '
$ws.Cells.Item(11, 4).Value = "'19.84"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -3.01%  '
$ws.Cells.Item(12, 5).Value = '  -1.44%  '
$ws.Cells.Item(13, 5).Value = '  +1.48%  '
$ws.Cells.Item(14, 4).Value = "'7.35"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.30%  '
$ws.Cells.Item(15, 4).Value = '3.048.44'
$ws.Cells.Item(15, 5).Value = '  -0.49%  '
$ws.Cells.Item(16, 4).Value = '2.630.43'
$ws.Cells.Item(16, 5).Value = '  -0.38%  '
$ws.Cells.Item(17, 4).Value = "'0.850"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -3.50%  '
$ws.Cells.Item(18, 4).Value = '49.484.70'
$ws.Cells.Item(18, 5).Value = '  -0.85%  '
$ws.Cells.Item(19, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(19, 4).Value = "'12.89"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.76%  '
$ws.Cells.Item(20, 2).Value = 'ImmutableX'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(20, 4).Value = "'2.95"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.57%  '
$ws.Cells.Item(21, 5).Value = '  -2.05%  '
$ws.Cells.Item(23, 4).Value = "'270.40"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -3.69%  '
$ws.Cells.Item(24, 4).Value = "'68.95"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -5.48%  '
$ws.Cells.Item(25, 5).Value = '  -2.37%  '
$ws.Cells.Item(26, 5).Value = '  -3.11%  '
$ws.Cells.Item(27, 4).Value = "'0.999"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.01%  '
$ws.Cells.Item(28, 4).Value = "'10.31"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +3.48%  '
$ws.Cells.Item(29, 5).Value = '  -1.27%  '
$ws.Cells.Item(30, 5).Value = '  -4.69%  '
$ws.Cells.Item(31, 4).Value = "'34.69"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -5.91%  '
$ws.Cells.Item(32, 4).Value = "'49.46"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.73%  '
$ws.Cells.Item(33, 5).Value = '  +0.65%  '
$ws.Cells.Item(34, 5).Value = '  +1.92%  '
$ws.Cells.Item(35, 5).Value = '  -0.16%  '
$ws.Cells.Item(36, 4).Value = "'19.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -3.30%  '
$ws.Cells.Item(37, 5).Value = '  +2.33%  '
$ws.Cells.Item(38, 5).Value = '  -1.33%  '
$ws.Cells.Item(39, 5).Value = '  -0.21%  '
$ws.Cells.Item(40, 4).Value = "'128.44"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.97%  '
$ws.Cells.Item(41, 5).Value = '  -1.80%  '
$ws.Cells.Item(42, 4).Value = "'22.15"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -1.91%  '
$ws.Cells.Item(43, 5).Value = '  +3.12%  '
$ws.Cells.Item(44, 4).Value = "'2.16"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -4.48%  '
$ws.Cells.Item(45, 4).Value = '2.058.35'
$ws.Cells.Item(45, 5).Value = '  -0.65%  '
$ws.Cells.Item(46, 2).Value = 'NEARProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(46, 4).Value = "'3.21"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -5.42%  '
$ws.Cells.Item(47, 2).Value = 'Stacks'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(47, 4).Value = "'2.12"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +6.53%  '
$ws.Cells.Item(48, 5).Value = '  -4.78%  '
$ws.Cells.Item(49, 4).Value = "'8.91"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -1.92%  '
$ws.Cells.Item(50, 4).Value = "'59.09"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +2.40%  '
$ws.Cells.Item(51, 5).Value = '  -3.73%  '
